$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '27.122.84'
Set-TextValue 'E2' '  +1.37%  '
Set-TextValue 'D3' '1.568.94'
Set-TextValue 'E3' '  +2.04%  '
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '208.25'
Set-TextValue 'E5' '  +1.28%  '
Set-TextValue 'E6' '  +1.19%  '
Set-TextValue 'E7' '  +0.18%  '
Set-TextValue 'D8' '22.14'
Set-TextValue 'E8' '  +4.23%  '
Set-TextValue 'E9' '  +1.37%  '
Set-TextValue 'D10' '0.0589'
Set-TextValue 'E10' '  +1.48%  '
Set-TextValue 'D11' '0.0861'
Set-TextValue 'E11' '  +0.82%  '
Set-TextValue 'D12' '1.790.97'
Set-TextValue 'D13' '1.567.87'
Set-TextValue 'E13' '  +1.78%  '
Set-TextValue 'E14' '  +2.78%  '
Set-TextValue 'D15' '0.523'
Set-TextValue 'E15' '  +2.82%  '
Set-TextValue 'D16' '27.108.24'
Set-TextValue 'E16' '  +1.29%  '
Set-TextValue 'D17' '62.09'
Set-TextValue 'E17' '  +1.89%  '
Set-TextValue 'D18' '219.54'
Set-TextValue 'E18' '  +3.00%  '
Set-TextValue 'D20' '7.37'
Set-TextValue 'E20' '  +1.72%  '
Set-TextValue 'E21' '  +0.23%  '
Set-TextValue 'E22' '  +2.03%  '
Set-TextValue 'E23' '  +1.72%  '
Set-TextValue 'E24' '  +1.61%  '
Set-TextValue 'D25' '154.27'
Set-TextValue 'E25' '  +1.57%  '
Set-TextValue 'D26' '6.64'
Set-TextValue 'E26' '  +0.93%  '
Set-TextValue 'D27' '15.02'
Set-TextValue 'E27' '  +1.52%  '
Set-TextValue 'E28' '  +0.21%  '
Set-TextValue 'E29' '  +1.75%  '
Set-TextValue 'E30' '  +3.17%  '
Set-TextValue 'E31' '  +0.38%  '
Set-TextValue 'E32' '  +0.94%  '
Set-TextValue 'D33' '1.450.44'
Set-TextValue 'E33' '  +6.24%  '
Set-TextValue 'E34' '  +4.85%  '
Set-TextValue 'E35' '  +4.40%  '
Set-TextValue 'D36' '0.966'
Set-TextValue 'E36' '  +0.77%  '
Set-TextValue 'E37' '  +0.90%  '
Set-TextValue 'D38' '0.0165'
Set-TextValue 'E38' '  +0.72%  '
Set-TextValue 'E39' '  +0.81%  '
Set-TextValue 'E40' '  +1.74%  '
Set-TextValue 'E41' '  +0.05%  '
Set-TextValue 'E42' '  +0.27%  '
Set-TextValue 'E43' '  +3.96%  '
Set-TextValue 'D45' '64.70'
Set-TextValue 'E45' '  +2.94%  '
Set-TextValue 'E46' '  +2.24%  '
Set-TextValue 'D47' '1.706.40'
Set-TextValue 'E47' '  +2.11%  '
Set-TextValue 'D48' '86.92'
Set-TextValue 'E48' '  +3.20%  '
Set-TextValue 'E49' '  +3.32%  '
Set-TextValue 'D50' '0.0₆0101'
Set-TextValue 'E50' '  +3.72%  '
Set-TextValue 'D51' '0.0967'
Set-TextValue 'E51' '  +2.55%  '
